# Auto-generated edit script for khl_probabilities_human_tour.xlsx
# Updates Summary (rows 2-5) and Cards_telegram (rows 2-5) sheets

$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item("Summary")
$wsCards = $wb.Worksheets.Item("Cards_telegram")

# ---- Summary row 2 ----
$wsSummary.Range("B2").Value = 45988.70833333334
$wsSummary.Range("C2").Value = "Металлург Мг"
$wsSummary.Range("D2").Value = "Авангард"
$wsSummary.Range("E2").Value = "Металлург Мг – Авангард"
$wsSummary.Range("F2").Value = 897821
$wsSummary.Range("G2").Value = "https://text.khl.ru/text/897821.html"
$wsSummary.Range("H2").Value = 6.3
$wsSummary.Range("I2").Value = 5.678571
$wsSummary.Range("J2").Value = 11.978571
$wsSummary.Range("K2").Value = 37.406257
$wsSummary.Range("L2").Value = 35.92968
$wsSummary.Range("M2").Value = 73.335938
$wsSummary.Range("N2").Value = 0.304904
$wsSummary.Range("O2").Value = 0.126149
$wsSummary.Range("P2").Value = 0.549013
$wsSummary.Range("Q2").Value = 3.279720830162936
$wsSummary.Range("R2").Value = 7.92713378623691
$wsSummary.Range("S2").Value = 1.821450493886301
$wsSummary.Range("T2").Value = 30.4904
$wsSummary.Range("U2").Value = 12.6149
$wsSummary.Range("V2").Value = 54.9013
$wsSummary.Range("W2").Value = 0.046563
$wsSummary.Range("X2").Value = 0.933504
$wsSummary.Range("Y2").Value = 1.071232688879748
$wsSummary.Range("Z2").Value = 0.100349
$wsSummary.Range("AA2").Value = 0.879717
$wsSummary.Range("AB2").Value = 1.136729198139856
$wsSummary.Range("AC2").Value = 0.183432
$wsSummary.Range("AD2").Value = 0.796634
$wsSummary.Range("AE2").Value = 1.255281597320727
$wsSummary.Range("AF2").Value = 0.917421
$wsSummary.Range("AG2").Value = 0.082579
$wsSummary.Range("AH2").Value = 1.090012110034542
$wsSummary.Range("AI2").Value = 0.780181
$wsSummary.Range("AJ2").Value = 0.219819
$wsSummary.Range("AK2").Value = 1.281753849427248
$wsSummary.Range("AL2").Value = 0.964012
$wsSummary.Range("AM2").Value = 0.035988
$wsSummary.Range("AN2").Value = 1.037331485500181
$wsSummary.Range("AO2").Value = 0.8866
$wsSummary.Range("AP2").Value = 0.1134
$wsSummary.Range("AQ2").Value = 1.127904353710805
$wsSummary.Range("AR2").Value = 0.563876
$wsSummary.Range("AS2").Value = 1.77343955053948
$wsSummary.Range("AT2").Value = 0.781871
$wsSummary.Range("AU2").Value = 1.278983361705448

# ---- Summary row 3 ----
$wsSummary.Range("B3").Value = 45988.70833333334
$wsSummary.Range("C3").Value = "Салават Юлаев"
$wsSummary.Range("D3").Value = "Барыс"
$wsSummary.Range("E3").Value = "Салават Юлаев – Барыс"
$wsSummary.Range("F3").Value = 897819
$wsSummary.Range("G3").Value = "https://text.khl.ru/text/897819.html"
$wsSummary.Range("H3").Value = 1.454545
$wsSummary.Range("I3").Value = 1
$wsSummary.Range("J3").Value = 2.454545
$wsSummary.Range("K3").Value = 23.48092
$wsSummary.Range("L3").Value = 23.645876
$wsSummary.Range("M3").Value = 47.126796
$wsSummary.Range("N3").Value = 0.59579
$wsSummary.Range("O3").Value = 0.172039
$wsSummary.Range("P3").Value = 0.231982
$wsSummary.Range("Q3").Value = 1.678443746957821
$wsSummary.Range("R3").Value = 5.812635507065258
$wsSummary.Range("S3").Value = 4.310679276840444
$wsSummary.Range("T3").Value = 59.57900000000001
$wsSummary.Range("U3").Value = 17.2039
$wsSummary.Range("V3").Value = 23.1982
$wsSummary.Range("W3").Value = 0.505953
$wsSummary.Range("X3").Value = 0.493858
$wsSummary.Range("Y3").Value = 2.024873546647012
$wsSummary.Range("Z3").Value = 0.679033
$wsSummary.Range("AA3").Value = 0.320779
$wsSummary.Range("AB3").Value = 3.117411052469146
$wsSummary.Range("AC3").Value = 0.81285
$wsSummary.Range("AD3").Value = 0.186961
$wsSummary.Range("AE3").Value = 5.348709089061355
$wsSummary.Range("AF3").Value = 0.775796
$wsSummary.Range("AG3").Value = 0.224204
$wsSummary.Range("AH3").Value = 1.288998654285405
$wsSummary.Range("AI3").Value = 0.540241
$wsSummary.Range("AJ3").Value = 0.459759
$wsSummary.Range("AK3").Value = 1.8510257459171
$wsSummary.Range("AL3").Value = 0.536572
$wsSummary.Range("AM3").Value = 0.463428
$wsSummary.Range("AN3").Value = 1.863682786280313
$wsSummary.Range("AO3").Value = 0.268847
$wsSummary.Range("AP3").Value = 0.731153
$wsSummary.Range("AQ3").Value = 3.719587720897016
$wsSummary.Range("AR3").Value = 0.888398
$wsSummary.Range("AS3").Value = 1.125621624542153
$wsSummary.Range("AT3").Value = 0.594522
$wsSummary.Range("AU3").Value = 1.682023541601488

# ---- Summary row 4 ----
$wsSummary.Range("B4").Value = 45988.8125
$wsSummary.Range("C4").Value = "Спартак"
$wsSummary.Range("D4").Value = "Автомобилист"
$wsSummary.Range("E4").Value = "Спартак – Автомобилист"
$wsSummary.Range("F4").Value = 897820
$wsSummary.Range("G4").Value = "https://text.khl.ru/text/897820.html"
$wsSummary.Range("H4").Value = 4.07223
$wsSummary.Range("I4").Value = 2.501555
$wsSummary.Range("J4").Value = 6.573785
$wsSummary.Range("K4").Value = 36.594732
$wsSummary.Range("L4").Value = 28.024318
$wsSummary.Range("M4").Value = 64.61905
$wsSummary.Range("N4").Value = 0.409069
$wsSummary.Range("O4").Value = 0.153703
$wsSummary.Range("P4").Value = 0.435191
$wsSummary.Range("Q4").Value = 2.444575365036216
$wsSummary.Range("R4").Value = 6.50605388313826
$wsSummary.Range("S4").Value = 2.297841637350037
$wsSummary.Range("T4").Value = 40.9069
$wsSummary.Range("U4").Value = 15.3703
$wsSummary.Range("V4").Value = 43.5191
$wsSummary.Range("W4").Value = 0.17336
$wsSummary.Range("X4").Value = 0.824603
$wsSummary.Range("Y4").Value = 1.212704780360974
$wsSummary.Range("Z4").Value = 0.301223
$wsSummary.Range("AA4").Value = 0.69674
$wsSummary.Range("AB4").Value = 1.435255619025748
$wsSummary.Range("AC4").Value = 0.450312
$wsSummary.Range("AD4").Value = 0.547651
$wsSummary.Range("AE4").Value = 1.825980414534074
$wsSummary.Range("AF4").Value = 0.859117
$wsSummary.Range("AG4").Value = 0.140883
$wsSummary.Range("AH4").Value = 1.163985813340907
$wsSummary.Range("AI4").Value = 0.670471
$wsSummary.Range("AJ4").Value = 0.329529
$wsSummary.Range("AK4").Value = 1.491488819054068
$wsSummary.Range("AL4").Value = 0.86853
$wsSummary.Range("AM4").Value = 0.13147
$wsSummary.Range("AN4").Value = 1.151370706826477
$wsSummary.Range("AO4").Value = 0.686932
$wsSummary.Range("AP4").Value = 0.313068
$wsSummary.Range("AQ4").Value = 1.455748167213057
$wsSummary.Range("AR4").Value = 0.706839
$wsSummary.Range("AS4").Value = 1.414749327640382
$wsSummary.Range("AT4").Value = 0.729345
$wsSummary.Range("AU4").Value = 1.371093241195867

# ---- Summary row 5 ----
$wsSummary.Range("B5").Value = 45988.8125
$wsSummary.Range("C5").Value = "ЦСКА"
$wsSummary.Range("D5").Value = "Лада"
$wsSummary.Range("E5").Value = "ЦСКА – Лада"
$wsSummary.Range("F5").Value = 897822
$wsSummary.Range("G5").Value = "https://text.khl.ru/text/897822.html"
$wsSummary.Range("H5").Value = 2.188722
$wsSummary.Range("I5").Value = 1.117647
$wsSummary.Range("J5").Value = 3.306369
$wsSummary.Range("K5").Value = 27.591019
$wsSummary.Range("L5").Value = 20.850574
$wsSummary.Range("M5").Value = 48.441593
$wsSummary.Range("N5").Value = 0.719666
$wsSummary.Range("O5").Value = 0.151933
$wsSummary.Range("P5").Value = 0.128229
$wsSummary.Range("Q5").Value = 1.389533478030086
$wsSummary.Range("R5").Value = 6.581848577991614
$wsSummary.Range("S5").Value = 7.798547910379087
$wsSummary.Range("T5").Value = 71.9666
$wsSummary.Range("U5").Value = 15.1933
$wsSummary.Range("V5").Value = 12.8229
$wsSummary.Range("W5").Value = 0.635914
$wsSummary.Range("X5").Value = 0.363914
$wsSummary.Range("Y5").Value = 2.747901976840682
$wsSummary.Range("Z5").Value = 0.790767
$wsSummary.Range("AA5").Value = 0.209062
$wsSummary.Range("AB5").Value = 4.783270034726541
$wsSummary.Range("AC5").Value = 0.893066
$wsSummary.Range("AD5").Value = 0.106762
$wsSummary.Range("AE5").Value = 9.36662857571046
$wsSummary.Range("AF5").Value = 0.771575
$wsSummary.Range("AG5").Value = 0.228425
$wsSummary.Range("AH5").Value = 1.296050286751126
$wsSummary.Range("AI5").Value = 0.534272
$wsSummary.Range("AJ5").Value = 0.465728
$wsSummary.Range("AK5").Value = 1.871705797795879
$wsSummary.Range("AL5").Value = 0.318542
$wsSummary.Range("AM5").Value = 0.681458
$wsSummary.Range("AN5").Value = 3.139303451350214
$wsSummary.Range("AO5").Value = 0.109459
$wsSummary.Range("AP5").Value = 0.890541
$wsSummary.Range("AQ5").Value = 9.135840817109603
$wsSummary.Range("AR5").Value = 0.953673
$wsSummary.Range("AS5").Value = 1.048577447405977
$wsSummary.Range("AT5").Value = 0.481443
$wsSummary.Range("AU5").Value = 2.077089084273736

# ---- Cards_telegram rows ----
# Cards_telegram row 2
$wsCards.Range("A2").Value = 45988.70833333334
$wsCards.Range("B2").Value = "Металлург Мг – Авангард"
$cardText2 = @"
КХЛ • Регулярный чемпионат • 27.11.2025

Металлург Мг – Авангард

Ожидания модели (60’):
• Голы: λ_total ≈ 9.27 (4.13 : 5.14)
• Броски: SOG λ ≈ 73 (37 : 36)

Исход (60’), честные кф:
• П1: 30.5%  (Kмод 3.28)
• Х:  12.6%  (Kмод 7.93)
• П2: 54.9%  (Kмод 1.82)

Тоталы голов:
• ТМ 4.5: 4.7%  (Kмод 21.48)
• ТБ 4.5: 93.4%  (Kмод 1.07)

• ТМ 5.5: 10.0%  (Kмод 9.97)
• ТБ 5.5: 88.0%  (Kмод 1.14)

• ТМ 6.5: 18.3%  (Kмод 5.45)
• ТБ 6.5: 79.7%  (Kмод 1.26)

Индивидуальные тоталы:
• Металлург Мг ИТБ 1.5: 91.7% (Kмод 1.09)
• Металлург Мг ИТБ 2.5: 78.0% (Kмод 1.28)
• Авангард ИТБ 1.5: 96.4% (Kмод 1.04)
• Авангард ИТБ 2.5: 88.7% (Kмод 1.13)

Фора +1.5:
• Металлург Мг +1.5: 56.4% (Kмод 1.77)
• Авангард +1.5: 78.2% (Kмод 1.28)
"@
$wsCards.Range("C2").Value = $cardText2

# Cards_telegram row 3
$wsCards.Range("A3").Value = 45988.70833333334
$wsCards.Range("B3").Value = "Салават Юлаев – Барыс"
$cardText3 = @"
КХЛ • Регулярный чемпионат • 27.11.2025

Салават Юлаев – Барыс

Ожидания модели (60’):
• Голы: λ_total ≈ 4.64 (2.84 : 1.80)
• Броски: SOG λ ≈ 47 (23 : 24)

Исход (60’), честные кф:
• П1: 59.6%  (Kмод 1.68)
• Х:  17.2%  (Kмод 5.81)
• П2: 23.2%  (Kмод 4.31)

Тоталы голов:
• ТМ 4.5: 50.6%  (Kмод 1.98)
• ТБ 4.5: 49.4%  (Kмод 2.02)

• ТМ 5.5: 67.9%  (Kмод 1.47)
• ТБ 5.5: 32.1%  (Kмод 3.12)

• ТМ 6.5: 81.3%  (Kмод 1.23)
• ТБ 6.5: 18.7%  (Kмод 5.35)

Индивидуальные тоталы:
• Салават Юлаев ИТБ 1.5: 77.6% (Kмод 1.29)
• Салават Юлаев ИТБ 2.5: 54.0% (Kмод 1.85)
• Барыс ИТБ 1.5: 53.7% (Kмод 1.86)
• Барыс ИТБ 2.5: 26.9% (Kмод 3.72)

Фора +1.5:
• Салават Юлаев +1.5: 88.8% (Kмод 1.13)
• Барыс +1.5: 59.5% (Kмод 1.68)
"@
$wsCards.Range("C3").Value = $cardText3

# Cards_telegram row 4
$wsCards.Range("A4").Value = 45988.8125
$wsCards.Range("B4").Value = "Спартак – Автомобилист"
$cardText4 = @"
КХЛ • Регулярный чемпионат • 27.11.2025

Спартак – Автомобилист

Ожидания модели (60’):
• Голы: λ_total ≈ 7.00 (3.45 : 3.54)
• Броски: SOG λ ≈ 65 (37 : 28)

Исход (60’), честные кф:
• П1: 40.9%  (Kмод 2.44)
• Х:  15.4%  (Kмод 6.51)
• П2: 43.5%  (Kмод 2.30)

Тоталы голов:
• ТМ 4.5: 17.3%  (Kмод 5.77)
• ТБ 4.5: 82.5%  (Kмод 1.21)

• ТМ 5.5: 30.1%  (Kмод 3.32)
• ТБ 5.5: 69.7%  (Kмод 1.44)

• ТМ 6.5: 45.0%  (Kмод 2.22)
• ТБ 6.5: 54.8%  (Kмод 1.83)

Индивидуальные тоталы:
• Спартак ИТБ 1.5: 85.9% (Kмод 1.16)
• Спартак ИТБ 2.5: 67.0% (Kмод 1.49)
• Автомобилист ИТБ 1.5: 86.9% (Kмод 1.15)
• Автомобилист ИТБ 2.5: 68.7% (Kмод 1.46)

Фора +1.5:
• Спартак +1.5: 70.7% (Kмод 1.41)
• Автомобилист +1.5: 72.9% (Kмод 1.37)
"@
$wsCards.Range("C4").Value = $cardText4

# Cards_telegram row 5
$wsCards.Range("A5").Value = 45988.8125
$wsCards.Range("B5").Value = "ЦСКА – Лада"
$cardText5 = @"
КХЛ • Регулярный чемпионат • 27.11.2025

ЦСКА – Лада

Ожидания модели (60’):
• Голы: λ_total ≈ 3.96 (2.82 : 1.15)
• Броски: SOG λ ≈ 48 (28 : 21)

Исход (60’), честные кф:
• П1: 72.0%  (Kмод 1.39)
• Х:  15.2%  (Kмод 6.58)
• П2: 12.8%  (Kмод 7.80)

Тоталы голов:
• ТМ 4.5: 63.6%  (Kмод 1.57)
• ТБ 4.5: 36.4%  (Kмод 2.75)

• ТМ 5.5: 79.1%  (Kмод 1.26)
• ТБ 5.5: 20.9%  (Kмод 4.78)

• ТМ 6.5: 89.3%  (Kмод 1.12)
• ТБ 6.5: 10.7%  (Kмод 9.37)

Индивидуальные тоталы:
• ЦСКА ИТБ 1.5: 77.2% (Kмод 1.30)
• ЦСКА ИТБ 2.5: 53.4% (Kмод 1.87)
• Лада ИТБ 1.5: 31.9% (Kмод 3.14)
• Лада ИТБ 2.5: 10.9% (Kмод 9.14)

Фора +1.5:
• ЦСКА +1.5: 95.4% (Kмод 1.05)
• Лада +1.5: 48.1% (Kмод 2.08)
"@
$wsCards.Range("C5").Value = $cardText5

